$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 21387.54
$ws.Range("I21").Value = 15912.546
$ws.Range("J21").Value = 51500
$ws.Range("K21").Value = 15912.546
$ws.Range("L21").Value = 51500
$ws.Range("M21").Value = -15444.546
$ws.Range("N21").Value = -52436

$ws.Range("H23").Value = 21387.54
$ws.Range("I23").Value = 15912.546
$ws.Range("J23").Value = 51500
$ws.Range("K23").Value = 15912.546
$ws.Range("L23").Value = 51500
$ws.Range("M23").Value = -15678.546
$ws.Range("N23").Value = -51968

$ws.Range("H131").Value = 1523.6666
$ws.Range("I131").Value = 1604.2307
$ws.Range("J131").Value = 1000
$ws.Range("K131").Value = 4812.6921
$ws.Range("L131").Value = 3000
$ws.Range("M131").Value = 227.3078999999998
$ws.Range("N131").Value = -13080

$ws.Range("H137").Value = 1921.6
$ws.Range("I137").Value = 1319.0667
$ws.Range("J137").Value = 2524.1333
$ws.Range("K137").Value = 3957.2001
$ws.Range("L137").Value = 7572.3999
$ws.Range("M137").Value = -1407.2001
$ws.Range("N137").Value = -12672.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5885231
$ws.Range("I2").Value = 4297
$ws.Range("J2").Value = 14706632
$ws.Range("K2").Value = 4297
$ws.Range("L2").Value = 14706632
$ws.Range("M2").Value = -4184
$ws.Range("N2").Value = -14706858

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H96").Value = 35000
$ws.Range("J96").Value = 35000
$ws.Range("L96").Value = 35000
$ws.Range("N96").Value = -40492

$ws.Range("H116").Value = 5885231
$ws.Range("I116").Value = 4297
$ws.Range("J116").Value = 14706632
$ws.Range("K116").Value = 4297
$ws.Range("L116").Value = 14706632
$ws.Range("M116").Value = -2003
$ws.Range("N116").Value = -14711220

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5885231
$ws.Range("I3").Value = 4297
$ws.Range("J3").Value = 14706632
$ws.Range("K3").Value = 4297
$ws.Range("L3").Value = 14706632
$ws.Range("M3").Value = -4183
$ws.Range("N3").Value = -14706860

$ws.Range("H64").Value = 1443.3704
$ws.Range("I64").Value = 3352.5715
$ws.Range("J64").Value = 775.15
$ws.Range("K64").Value = 3352.5715
$ws.Range("L64").Value = 775.15
$ws.Range("M64").Value = -3127.5715
$ws.Range("N64").Value = -1225.15

$ws.Range("H67").Value = 1443.3704
$ws.Range("I67").Value = 3352.5715
$ws.Range("J67").Value = 775.15
$ws.Range("K67").Value = 3352.5715
$ws.Range("L67").Value = 775.15
$ws.Range("M67").Value = -2572.5715
$ws.Range("N67").Value = -2335.15

$ws.Range("H95").Value = 16750
$ws.Range("J95").Value = 16750
$ws.Range("L95").Value = 16750
$ws.Range("N95").Value = -22242

$ws.Range("H134").Value = 25714.6
$ws.Range("I134").Value = 26230.842
$ws.Range("K134").Value = 78692.526
$ws.Range("M134").Value = -76157.526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 20000
$ws.Range("J68").Value = 20000
$ws.Range("L68").Value = 20000
$ws.Range("N68").Value = -21498

$ws.Range("H71").Value = 20000
$ws.Range("J71").Value = 20000
$ws.Range("L71").Value = 60000
$ws.Range("N71").Value = -67488

$ws.Range("H92").Value = 29314.428
$ws.Range("J92").Value = 29314.428
$ws.Range("L92").Value = 29314.428
$ws.Range("N92").Value = -34306.428

$ws.Range("H96").Value = 16571.143
$ws.Range("J96").Value = 16571.143
$ws.Range("L96").Value = 16571.143
$ws.Range("N96").Value = -22063.143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 25001402
$ws.Range("J34").Value = 29413380
$ws.Range("L34").Value = 88240140
$ws.Range("N34").Value = -88240308

$ws.Range("H49").Value = 3550
$ws.Range("J49").Value = 3550
$ws.Range("L49").Value = 10650
$ws.Range("N49").Value = -10962

$ws.Range("H100").Value = 1980
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H104").Value = 3300
$ws.Range("J104").Value = 3300
$ws.Range("L104").Value = 9900
$ws.Range("N104").Value = -15142

$ws.Range("H131").Value = 1925827.2
$ws.Range("J131").Value = 2175881.5
$ws.Range("L131").Value = 6527644.5
$ws.Range("N131").Value = -6537724.5

$ws.Range("H134").Value = 3829.2559
$ws.Range("I134").Value = 1886.8462
$ws.Range("J134").Value = 6800
$ws.Range("K134").Value = 5660.5386
$ws.Range("L134").Value = 20400
$ws.Range("M134").Value = -590.5385999999999
$ws.Range("N134").Value = -30540

$ws.Range("H139").Value = 1464.95
$ws.Range("I139").Value = 1081.1875
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 3243.5625
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 1896.4375
$ws.Range("N139").Value = -19280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1775.2858
$ws.Range("I113").Value = 2028
$ws.Range("J113").Value = 1438.3334
$ws.Range("K113").Value = 2028
$ws.Range("L113").Value = 1438.3334
$ws.Range("M113").Value = 142
$ws.Range("N113").Value = -5778.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 12999.333
$ws.Range("I26").Value = 14000
$ws.Range("J26").Value = 12499
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 12499
$ws.Range("M26").Value = -13705
$ws.Range("N26").Value = -13089

$ws.Range("H55").Value = 178.90475
$ws.Range("I55").Value = 165.92308
$ws.Range("J55").Value = 200
$ws.Range("K55").Value = 165.92308
$ws.Range("L55").Value = 200
$ws.Range("M55").Value = 7.076920000000001
$ws.Range("N55").Value = -546

$ws.Range("H61").Value = 949.625
$ws.Range("I61").Value = 942.4286
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 942.4286
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -740.4286
$ws.Range("N61").Value = -1404

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H104").Value = 17166.666
$ws.Range("J104").Value = 17166.666
$ws.Range("L104").Value = 17166.666
$ws.Range("N104").Value = -24154.666

$ws.Range("H113").Value = 949.625
$ws.Range("I113").Value = 942.4286
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 942.4286
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1227.5714
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 70017
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H35").Value = 70017
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

$ws.Range("H98").Value = 5000
$ws.Range("J98").Value = 5000
$ws.Range("L98").Value = 5000
$ws.Range("N98").Value = -10990

$ws.Range("H104").Value = 22561.5
$ws.Range("J104").Value = 22561.5
$ws.Range("L104").Value = 22561.5
$ws.Range("N104").Value = -29549.5

$ws.Range("H105").Value = 27400
$ws.Range("J105").Value = 27400
$ws.Range("L105").Value = 27400
$ws.Range("N105").Value = -34388

$ws.Range("H136").Value = 887.43243
$ws.Range("I136").Value = 818.6923
$ws.Range("K136").Value = 2456.0769
$ws.Range("M136").Value = 93.92309999999998
